# Update cryptos list - GitHub Actions scheduled refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "63.486.83"
$ws.Range("E2").Value = "  +5.72%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.391.19"
$ws.Range("E3").Value = "  +6.14%  "

# Row 5 - BNB
$ws.Range("D5").Value = "576.05"
$ws.Range("E5").Value = "  +7.30%  "

# Row 6 - Solana
$ws.Range("D6").Value = "154.51"
$ws.Range("E6").Value = "  +6.20%  "

# Row 7 - USDC
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.14%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.394.60"
$ws.Range("E8").Value = "  +6.04%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  +0.09%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  +1.95%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "0.121"
$ws.Range("E11").Value = "  +7.14%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +1.55%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.972.52"
$ws.Range("E13").Value = "  +6.04%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +0.39%  "

# Row 15 - ShibaInu
$ws.Range("D15").Value = "0.0000185"
$ws.Range("E15").Value = "  +7.05%  "

# Row 16 - Avalanche
$ws.Range("D16").Value = "27.17"
$ws.Range("E16").Value = "  +5.01%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "63.501.84"
$ws.Range("E17").Value = "  +5.74%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.396.70"
$ws.Range("E18").Value = "  +6.43%  "

# Row 19 - Polkadot
$ws.Range("D19").Value = "6.37"
$ws.Range("E19").Value = "  +1.44%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "13.95"
$ws.Range("E20").Value = "  +4.71%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "8.45"
$ws.Range("E21").Value = "  +2.90%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "389.13"
$ws.Range("E22").Value = "  +5.29%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.17%  "

# Row 24 - Polygon
$ws.Range("D24").Value = "0.537"
$ws.Range("E24").Value = "  +2.76%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "70.93"
$ws.Range("E25").Value = "  +2.17%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("D26").Value = "9.74"
$ws.Range("E26").Value = "  +12.47%  "

# Row 27 - PEPE
$ws.Range("D27").Value = "0.0000104"
$ws.Range("E27").Value = "  +18.05%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  +6.32%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.12%  "

# Row 30 - PancakeSwap
$ws.Range("D30").Value = "2.03"
$ws.Range("E30").Value = "  +7.31%  "

# Row 31 - RenderToken
$ws.Range("D31").Value = "6.45"
$ws.Range("E31").Value = "  +5.06%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "23.16"
$ws.Range("E32").Value = "  +2.93%  "

# Row 33 - now Fetch.AI (was NEARProtocol) - rows 33/34 swapped content
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "1.32"
$ws.Range("E33").Value = "  +10.34%  "

# Row 34 - now NEARProtocol (was Fetch.AI)
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "5.59"
$ws.Range("E34").Value = "  +6.04%  "

# Row 35 - Aptos
$ws.Range("D35").Value = "6.76"
$ws.Range("E35").Value = "  +2.69%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "1.48"
$ws.Range("E36").Value = "  +8.87%  "

# Row 37 - Monero
$ws.Range("D37").Value = "158.28"
$ws.Range("E37").Value = "  +1.45%  "

# Row 38 - EnergySwap
$ws.Range("D38").Value = "27.85"
$ws.Range("E38").Value = "  +5.79%  "

# Row 39 - Stacks
$ws.Range("E39").Value = "  +12.48%  "

# Row 40 - Maker
$ws.Range("D40").Value = "2.913.92"
$ws.Range("E40").Value = "  +2.49%  "

# Row 41 - Hedera
$ws.Range("D41").Value = "0.0746"
$ws.Range("E41").Value = "  +5.56%  "

# Row 42 - VeChain
$ws.Range("E42").Value = "  +6.47%  "

# Row 43 - Mantle
$ws.Range("E43").Value = "  +5.85%  "

# Row 44 - now OKB (was Filecoin) - rows 44/45 swapped content
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "41.11"
$ws.Range("E44").Value = "  +2.95%  "

# Row 45 - now Filecoin (was OKB)
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").Value = "4.31"
$ws.Range("E45").Value = "  +1.48%  "

# Row 46 - ONDO
$ws.Range("E46").Value = "  +7.14%  "

# Row 47 - RenzoRestakedETH
$ws.Range("D47").Value = "3.435.01"
$ws.Range("E47").Value = "  +6.10%  "

# Row 48 - InjectiveProtocol
$ws.Range("D48").Value = "22.11"
$ws.Range("E48").Value = "  +6.46%  "

# Row 49 - Bittensor
$ws.Range("D49").Value = "300.92"
$ws.Range("E49").Value = "  +13.56%  "

# Row 50 - Stellar
$ws.Range("E50").Value = "  -1.00%  "

# Row 51 - Cosmos
$ws.Range("E51").Value = "  +2.53%  "
